$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (price + 1h volume change columns).
# A leading apostrophe forces Excel to store a numeric-looking Price value
# (e.g. "1.00") as literal text, matching the source data's text formatting
# instead of letting AutoDetect coerce it into a Number.
$ws.Cells.Item(2, 4).Value = '47.347.26'
$ws.Cells.Item(2, 5).Value = '  +0.23%  '
$ws.Cells.Item(3, 4).Value = '2.488.34'
$ws.Cells.Item(3, 5).Value = '  -0.26%  '
$ws.Cells.Item(4, 4).Value = '''1.00'
$ws.Cells.Item(4, 5).Value = '  +0.13%  '
$ws.Cells.Item(5, 4).Value = '''321.15'
$ws.Cells.Item(5, 5).Value = '  -0.39%  '
$ws.Cells.Item(6, 4).Value = '''108.34'
$ws.Cells.Item(6, 5).Value = '  +2.41%  '
$ws.Cells.Item(7, 5).Value = '  +0.00%  '
$ws.Cells.Item(8, 4).Value = '''1.00'
$ws.Cells.Item(8, 5).Value = '  +0.10%  '
$ws.Cells.Item(9, 4).Value = '''0.536'
$ws.Cells.Item(9, 5).Value = '  -0.78%  '
$ws.Cells.Item(10, 4).Value = '''39.21'
$ws.Cells.Item(10, 5).Value = '  +4.56%  '
$ws.Cells.Item(11, 5).Value = '  -0.54%  '
$ws.Cells.Item(12, 5).Value = '  +0.43%  '
$ws.Cells.Item(13, 5).Value = '  +0.14%  '
$ws.Cells.Item(14, 4).Value = '''7.13'
$ws.Cells.Item(14, 5).Value = '  -0.56%  '
$ws.Cells.Item(15, 4).Value = '2.877.99'
$ws.Cells.Item(15, 5).Value = '  -0.19%  '
$ws.Cells.Item(16, 4).Value = '2.491.01'
$ws.Cells.Item(16, 5).Value = '  -1.23%  '
$ws.Cells.Item(17, 4).Value = '''0.843'
$ws.Cells.Item(17, 5).Value = '  -0.10%  '
$ws.Cells.Item(18, 4).Value = '47.247.73'
$ws.Cells.Item(18, 5).Value = '  +0.19%  '
$ws.Cells.Item(19, 4).Value = '''13.13'
$ws.Cells.Item(19, 5).Value = '  +3.07%  '
$ws.Cells.Item(20, 4).Value = '''6.63'
$ws.Cells.Item(20, 5).Value = '  +1.24%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0935'
$ws.Cells.Item(21, 5).Value = '  +0.07%  '
$ws.Cells.Item(22, 4).Value = '''2.67'
$ws.Cells.Item(22, 5).Value = '  +13.04%  '
$ws.Cells.Item(23, 4).Value = '''70.37'
$ws.Cells.Item(23, 5).Value = '  -0.67%  '
$ws.Cells.Item(24, 4).Value = '''245.06'
$ws.Cells.Item(24, 5).Value = '  -2.45%  '
$ws.Cells.Item(25, 4).Value = '''2.55'
$ws.Cells.Item(25, 5).Value = '  +0.51%  '
$ws.Cells.Item(26, 4).Value = '''1.00'
$ws.Cells.Item(26, 5).Value = '  -0.01%  '
$ws.Cells.Item(27, 4).Value = '''25.71'
$ws.Cells.Item(27, 5).Value = '  -1.69%  '
$ws.Cells.Item(28, 4).Value = '''2.27'
$ws.Cells.Item(28, 5).Value = '  +3.52%  '
$ws.Cells.Item(29, 4).Value = '''9.97'
$ws.Cells.Item(29, 5).Value = '  -1.95%  '
$ws.Cells.Item(30, 5).Value = '  +3.09%  '
$ws.Cells.Item(31, 4).Value = '''34.53'
$ws.Cells.Item(31, 5).Value = '  -2.28%  '
$ws.Cells.Item(32, 4).Value = '''49.79'
$ws.Cells.Item(32, 5).Value = '  +0.37%  '
$ws.Cells.Item(33, 4).Value = '''20.51'
$ws.Cells.Item(33, 5).Value = '  +3.45%  '
$ws.Cells.Item(34, 4).Value = '''5.34'
$ws.Cells.Item(34, 5).Value = '  -0.80%  '
$ws.Cells.Item(35, 4).Value = '''0.0785'
$ws.Cells.Item(35, 5).Value = '  +0.28%  '
$ws.Cells.Item(36, 4).Value = '''1.00'
$ws.Cells.Item(36, 5).Value = '  +0.22%  '
$ws.Cells.Item(37, 4).Value = '''4.73'
$ws.Cells.Item(37, 5).Value = '  +2.08%  '
$ws.Cells.Item(38, 4).Value = '''1.96'
$ws.Cells.Item(38, 5).Value = '  +1.29%  '
$ws.Cells.Item(39, 5).Value = '  -1.88%  '
$ws.Cells.Item(40, 4).Value = '''23.33'
$ws.Cells.Item(40, 5).Value = '  +6.73%  '
$ws.Cells.Item(41, 5).Value = '  -0.11%  '
$ws.Cells.Item(43, 4).Value = '''117.24'
$ws.Cells.Item(43, 5).Value = '  -3.42%  '
$ws.Cells.Item(44, 5).Value = '  +0.42%  '
$ws.Cells.Item(45, 4).Value = '1.997.62'
$ws.Cells.Item(45, 5).Value = '  +2.30%  '
$ws.Cells.Item(46, 4).Value = '''3.03'
$ws.Cells.Item(46, 5).Value = '  +1.82%  '
$ws.Cells.Item(47, 5).Value = '  -5.17%  '
$ws.Cells.Item(48, 4).Value = '''9.15'
$ws.Cells.Item(48, 5).Value = '  -0.33%  '
$ws.Cells.Item(49, 4).Value = '''1.78'
$ws.Cells.Item(49, 5).Value = '  -1.18%  '
$ws.Cells.Item(50, 5).Value = '  -5.61%  '
$ws.Cells.Item(51, 4).Value = '''56.37'
$ws.Cells.Item(51, 5).Value = '  +2.58%  '
